$wb = $excel.ActiveWorkbook

$sheetNames = @("Step3_DataPts_0.5", "Step3_DataPts_0.7", "Step3_DataPts_0.8", "Step3_DataPts_0.9")

# Column G (Pulse_Width) values differ per sheet/pressure, keyed by row.
$gValues = @{
    2 = @(16, 55, 65, 77)
    3 = @(18, 58, 68, 85)
    4 = @(72, 97, 111, 127)
    5 = @(23, 59, 70, 89)
    6 = @(18, 55, 64, 77)
}

# Column C (First_Noticeable_Increase_Index) and E (First_Noticeable_Increase_Cumulative_Value)
# are identical across all four pressure sheets for a given row.
$cValues = @{
    2 = 87
    3 = 88
    4 = 35
    5 = 87
    6 = 88
}

$eValues = @{
    2 = 0.006731925240992658
    3 = 0.01651608338647059
    4 = 0.06015755477446096
    5 = 0.01826365714201592
    6 = 0.01195572242935386
}

for ($i = 0; $i -lt $sheetNames.Count; $i++) {
    $ws = $wb.Worksheets.Item($sheetNames[$i])
    foreach ($row in 2..6) {
        $ws.Cells.Item($row, 3).Value = $cValues[$row]
        $ws.Cells.Item($row, 5).Value = $eValues[$row]
        $ws.Cells.Item($row, 7).Value = $gValues[$row][$i]
    }
}
